# fix(publipostage): Correct status name
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# statut_label: "bleu" -> "noir" (row 2)
$ws.Range("B2").Value = "noir"

# statut_name: correct wording for rows 2-4
$ws.Range("C2").Value = "pas de résultat postés ni publiés"
$ws.Range("C3").Value = "résultat postés ou publiés dans les 12 mois"
$ws.Range("C4").Value = "résultat postés ou publiés dans les 12 mois"
